# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.695.49'
$ws.Range('E2').Value = '  +5.86%  '
$ws.Range('D3').Value = '3.113.68'
$ws.Range('E3').Value = '  +3.83%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''586.03'
$ws.Range('E5').Value = '  +4.37%  '
$ws.Range('D6').Value = '''143.57'
$ws.Range('E6').Value = '  +3.85%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.101.28'
$ws.Range('E8').Value = '  +3.56%  '
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('D10').Value = '''0.145'
$ws.Range('E10').Value = '  +9.17%  '
$ws.Range('D11').Value = '''5.74'
$ws.Range('E11').Value = '  +9.68%  '
$ws.Range('D12').Value = '''0.469'
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('D13').Value = '''0.0000244'
$ws.Range('E13').Value = '  +5.49%  '
$ws.Range('D14').Value = '''35.60'
$ws.Range('E14').Value = '  +5.66%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').Value = '''7.30'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').Value = '3.106.99'
$ws.Range('E17').Value = '  +3.46%  '
$ws.Range('D18').Value = '62.681.31'
$ws.Range('E18').Value = '  +5.68%  '
$ws.Range('D19').Value = '''454.26'
$ws.Range('E19').Value = '  +5.63%  '
$ws.Range('D20').Value = '''14.10'
$ws.Range('E20').Value = '  +3.38%  '
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('E22').Value = '  +5.83%  '
$ws.Range('E23').Value = '  +1.91%  '
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').Value = '''2.70'
$ws.Range('E27').Value = '  +5.79%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +5.10%  '
$ws.Range('D30').Value = '''6.86'
$ws.Range('E30').Value = '  +12.92%  '
$ws.Range('E31').Value = '  +12.50%  '
$ws.Range('D32').Value = '''27.13'
$ws.Range('E32').Value = '  +5.12%  '
$ws.Range('E33').Value = '  +5.11%  '
$ws.Range('D34').Value = '0.0₃0804'
$ws.Range('E34').Value = '  +5.44%  '
$ws.Range('D35').Value = '''6.08'
$ws.Range('E35').Value = '  +1.91%  '
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('E37').Value = '  +3.53%  '
$ws.Range('E38').Value = '  +9.96%  '
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').Value = '''425.02'
$ws.Range('E40').Value = '  +5.35%  '
$ws.Range('D41').Value = '2.949.09'
$ws.Range('E41').Value = '  +6.57%  '
$ws.Range('D42').Value = '''0.0373'
$ws.Range('E42').Value = '  +5.35%  '
$ws.Range('D43').Value = '''0.281'
$ws.Range('E43').Value = '  +11.19%  '
$ws.Range('E44').Value = '  +3.23%  '
$ws.Range('E45').Value = '  +7.95%  '
$ws.Range('D46').Value = '''125.32'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''35.19'
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '''0.999'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').Value = '''24.89'
$ws.Range('E50').Value = '  +5.79%  '
$ws.Range('E51').Value = '  +6.63%  '
